# Component three / C5-PowerPoint.pptx
#
# The table on the "SOURCES OF FINANCE" slide (slide 6) had its table
# style changed from the deck's custom "Table_0" style
# ({4AA9416C-67E9-481D-B133-AFE43CCA5107}) to the built-in
# "Medium Style 2 - Accent 1" table style
# ({B47683D6-7684-4547-AA33-42CA8073DDF9}).
#
# This mirrors picking a different style from the Table Styles gallery
# on the (Table Design) ribbon while the table is selected.

$newStyleId = "{B47683D6-7684-4547-AA33-42CA8073DDF9}"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# The table lives in the 2nd shape on the slide (the 1st is the title
# textbox); it's hosted in a graphic frame, so grab its .Table. Walk the
# shapes defensively in case ordering ever differs.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle($newStyleId)
    }
}
